$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 13; this shifts the existing rows 13-37
# (and all their formatting/content) down to rows 14-38, matching the
# diff which shows every record from the old row 13 onward moving down
# by one row, with a brand-new record appearing at the top (row 13).
$ws.Rows(13).Insert()

# Populate the newly inserted row 13 with the new record's data.
$ws.Cells.Item(13, 1).Value = 6
$ws.Cells.Item(13, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(13, 3).Value = "Metropolitana"
$ws.Cells.Item(13, 4).Value = 44690
$ws.Cells.Item(13, 5).Value = 13
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100101
$ws.Cells.Item(13, 8).Value = "Berries"
$ws.Cells.Item(13, 9).Value = 100101006
$ws.Cells.Item(13, 10).Value = "Higo"
$ws.Cells.Item(13, 11).Value = "Sin especificar"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 100
$ws.Cells.Item(13, 14).Value = 12000
$ws.Cells.Item(13, 15).Value = 12000
$ws.Cells.Item(13, 16).Value = 12000
$ws.Cells.Item(13, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(13, 18).Value = "Región Metropolitana"
$ws.Cells.Item(13, 19).Value = 1714
$ws.Cells.Item(13, 20).Value = 7
